$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename product description and bump cost (Bulto -> Tarima)
$ws.Range("C2").Value = "Tarima Charola 8x55 50 pzas 10 paquetes"
$ws.Range("E2").Value = 4000

# Row 3: rename product description (Bulto -> Tarima) and bump cost
$ws.Range("C3").Value = "Tarima Charola 8x55 50 pzas 10 paquetes Great Value"
$ws.Range("E3").Value = 4000

# Row 4: new product entry (CH857 / Oxxo) - copy formatting from row 2 first
$ws.Range("A2:F2").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)
$ws.Range("A4").Value = "CH857"
$ws.Range("B4").Value = 751095331
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 4000
$ws.Range("F4").Value = "P18"
$ws.Range("C4").Value = "Tarima Charola 8x55 50 pzas 10 paquetes Oxxo"

# Match the author's final active-cell selection
$ws.Range("G6").Select()
